$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (BF) held values like "4-30-2013-14" which were one day
# off from the actual game date because of how NBA stats were shown. Fix the
# training data by rewriting the BF column (rows 2-31) as "2014-04-30".
#
# Note: assigning a plain ISO-looking string via .Value would be
# auto-converted to a date serial by Excel's type inference, so we force a
# literal-text entry (leading apostrophe) and then reset the cell style back
# to "Normal" so no stray number-format/quote-prefix style sticks to the
# cell (matches the original formatting, which had no explicit style).
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)
    $cell.Value = "'2014-04-30"
    $cell.Style = "Normal"
}
